$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 143
$ws.Range("B143").Value = 6937250
$ws.Range("F143").Value = "Giannina"
$ws.Range("G143").Value = "Lamia"
$ws.Range("H143").Value = 1
$ws.Range("I143").Value = 4
$ws.Range("J143").Value = "A"
$ws.Range("K143").Value = 2.3
$ws.Range("L143").Value = 3.25
$ws.Range("M143").Value = 3.25
$ws.Range("N143").Value = 2.55
$ws.Range("O143").Value = 2.875
$ws.Range("P143").Value = 3.1
$ws.Range("Q143").Value = 0
$ws.Range("R143").Value = 1.75
$ws.Range("S143").Value = 2.125
$ws.Range("T143").Value = 2
$ws.Range("U143").Value = 1.85
$ws.Range("V143").Value = 2
$ws.Range("W143").Value = -1
$ws.Range("Y143").Value = 2.1
$ws.Range("Z143").Value = -1
$ws.Range("AA143").Value = 1.125
$ws.Range("AB143").Value = 0.8500000000000001

# Row 144
$ws.Range("B144").Value = 6937247
$ws.Range("F144").Value = "AEK Athens"
$ws.Range("G144").Value = "Asteras Tripolis"
$ws.Range("H144").Value = 4
$ws.Range("I144").Value = 2
$ws.Range("J144").Value = "H"
$ws.Range("K144").Value = 1.285
$ws.Range("L144").Value = 5.5
$ws.Range("M144").Value = 12
$ws.Range("N144").Value = 1.285
$ws.Range("O144").Value = 5.75
$ws.Range("P144").Value = 10
$ws.Range("Q144").Value = -1.5
$ws.Range("R144").Value = 1.825
$ws.Range("S144").Value = 2.025
$ws.Range("T144").Value = 3
$ws.Range("U144").Value = 2.025
$ws.Range("V144").Value = 1.825
$ws.Range("W144").Value = 0.2849999999999999
$ws.Range("Y144").Value = -1
$ws.Range("Z144").Value = 0.825
$ws.Range("AA144").Value = -1
$ws.Range("AB144").Value = 1.025

# Row 170
$ws.Range("B170").Value = 6937266
$ws.Range("F170").Value = "Atromitos Athinon"
$ws.Range("G170").Value = "Lamia"
$ws.Range("H170").Value = 3
$ws.Range("I170").Value = 1
$ws.Range("J170").Value = "H"
$ws.Range("K170").Value = 2.3
$ws.Range("L170").Value = 3.2
$ws.Range("M170").Value = 3.1
$ws.Range("N170").Value = 2.2
$ws.Range("O170").Value = 3.3
$ws.Range("P170").Value = 3.3
$ws.Range("Q170").Value = -0.25
$ws.Range("R170").Value = 1.925
$ws.Range("S170").Value = 1.925
$ws.Range("T170").Value = 2.5
$ws.Range("U170").Value = 2.025
$ws.Range("V170").Value = 1.825
$ws.Range("W170").Value = 1.2
$ws.Range("Y170").Value = -1
$ws.Range("Z170").Value = 0.925
$ws.Range("AA170").Value = -1
$ws.Range("AB170").Value = 1.025
$ws.Range("AC170").Value = -1

# Row 171
$ws.Range("B171").Value = 6937268
$ws.Range("F171").Value = "Panetolikos"
$ws.Range("G171").Value = "Olympiakos"
$ws.Range("H171").Value = 1
$ws.Range("I171").Value = 2
$ws.Range("J171").Value = "A"
$ws.Range("K171").Value = 8
$ws.Range("L171").Value = 5
$ws.Range("M171").Value = 1.363
$ws.Range("N171").Value = 8.5
$ws.Range("O171").Value = 5
$ws.Range("P171").Value = 1.363
$ws.Range("Q171").Value = 1.25
$ws.Range("R171").Value = 2.025
$ws.Range("S171").Value = 1.825
$ws.Range("T171").Value = 2.75
$ws.Range("U171").Value = 1.85
$ws.Range("V171").Value = 2
$ws.Range("W171").Value = -1
$ws.Range("Y171").Value = 0.363
$ws.Range("Z171").Value = 0.5125
$ws.Range("AA171").Value = -0.5
$ws.Range("AB171").Value = 0.425
$ws.Range("AC171").Value = -0.5

# Row 175
$ws.Range("B175").Value = 6937272
$ws.Range("F175").Value = "Lamia"
$ws.Range("G175").Value = "PAOK Salonika"
$ws.Range("H175").Value = 0
$ws.Range("I175").Value = 2
$ws.Range("J175").Value = "A"
$ws.Range("K175").Value = 7.5
$ws.Range("L175").Value = 4.5
$ws.Range("M175").Value = 1.444
$ws.Range("N175").Value = 9.5
$ws.Range("O175").Value = 5
$ws.Range("P175").Value = 1.333
$ws.Range("Q175").Value = 1.5
$ws.Range("R175").Value = 1.925
$ws.Range("S175").Value = 1.925
$ws.Range("T175").Value = 3
$ws.Range("U175").Value = 1.95
$ws.Range("V175").Value = 1.9
$ws.Range("W175").Value = -1
$ws.Range("Y175").Value = 0.333
$ws.Range("Z175").Value = -1
$ws.Range("AA175").Value = 0.925
$ws.Range("AB175").Value = -1
$ws.Range("AC175").Value = 0.8999999999999999

# Row 176
$ws.Range("B176").Value = 6935701
$ws.Range("F176").Value = "Kifisias FC"
$ws.Range("G176").Value = "Panetolikos"
$ws.Range("H176").Value = 2
$ws.Range("I176").Value = 2
$ws.Range("K176").Value = 2.45
$ws.Range("L176").Value = 3.25
$ws.Range("M176").Value = 3
$ws.Range("N176").Value = 2.05
$ws.Range("O176").Value = 3.3
$ws.Range("P176").Value = 3.8
$ws.Range("Q176").Value = -0.5
$ws.Range("T176").Value = 2.25
$ws.Range("U176").Value = 1.8
$ws.Range("V176").Value = 2.05
$ws.Range("X176").Value = 2.3
$ws.Range("Z176").Value = -1
$ws.Range("AA176").Value = 0.8
$ws.Range("AB176").Value = 0.8

# Row 177
$ws.Range("B177").Value = 6937271
$ws.Range("F177").Value = "Giannina"
$ws.Range("G177").Value = "Atromitos Athinon"
$ws.Range("H177").Value = 1
$ws.Range("I177").Value = 1
$ws.Range("K177").Value = 2.45
$ws.Range("L177").Value = 3.1
$ws.Range("M177").Value = 3.1
$ws.Range("N177").Value = 2
$ws.Range("O177").Value = 3.3
$ws.Range("P177").Value = 4
$ws.Range("Q177").Value = -0.5
$ws.Range("R177").Value = 2.025
$ws.Range("S177").Value = 1.825
$ws.Range("T177").Value = 2.25
$ws.Range("X177").Value = 2.3
$ws.Range("Z177").Value = -1
$ws.Range("AA177").Value = 0.825
$ws.Range("AB177").Value = -0.5
$ws.Range("AC177").Value = 0.5

# Row 178
$ws.Range("B178").Value = 6935700
$ws.Range("F178").Value = "Panserraikos"
$ws.Range("G178").Value = "Asteras Tripolis"
$ws.Range("I178").Value = 1
$ws.Range("J178").Value = "H"
$ws.Range("K178").Value = 2.6
$ws.Range("L178").Value = 3.2
$ws.Range("M178").Value = 2.875
$ws.Range("N178").Value = 2.25
$ws.Range("P178").Value = 3.3
$ws.Range("Q178").Value = -0.25
$ws.Range("R178").Value = 1.925
$ws.Range("S178").Value = 1.925
$ws.Range("U178").Value = 2
$ws.Range("V178").Value = 1.85
$ws.Range("W178").Value = 1.25
$ws.Range("X178").Value = -1
$ws.Range("Z178").Value = 0.925
$ws.Range("AA178").Value = -1
$ws.Range("AB178").Value = 1

# Row 179
$ws.Range("B179").Value = 6936863
$ws.Range("F179").Value = "OFI Crete"
$ws.Range("G179").Value = "Panathinaikos"
$ws.Range("I179").Value = 2
$ws.Range("J179").Value = "D"
$ws.Range("K179").Value = 8
$ws.Range("L179").Value = 4.75
$ws.Range("M179").Value = 1.4
$ws.Range("N179").Value = 5.5
$ws.Range("O179").Value = 4.75
$ws.Range("P179").Value = 1.55
$ws.Range("Q179").Value = 1
$ws.Range("R179").Value = 1.95
$ws.Range("S179").Value = 1.9
$ws.Range("T179").Value = 2.5
$ws.Range("U179").Value = 1.85
$ws.Range("V179").Value = 2
$ws.Range("W179").Value = -1
$ws.Range("X179").Value = 3.75
$ws.Range("Z179").Value = 0.95
$ws.Range("AB179").Value = 0.8500000000000001

# Row 180
$ws.Range("B180").Value = 6937270
$ws.Range("F180").Value = "Olympiakos"
$ws.Range("G180").Value = "Volos NFC"
$ws.Range("H180").Value = 3
$ws.Range("I180").Value = 0
$ws.Range("J180").Value = "H"
$ws.Range("K180").Value = 1.125
$ws.Range("L180").Value = 9
$ws.Range("M180").Value = 19
$ws.Range("N180").Value = 1.111
$ws.Range("O180").Value = 9
$ws.Range("P180").Value = 21
$ws.Range("Q180").Value = -2.25
$ws.Range("R180").Value = 1.875
$ws.Range("S180").Value = 1.975
$ws.Range("T180").Value = 3.25
$ws.Range("U180").Value = 2
$ws.Range("V180").Value = 1.85
$ws.Range("W180").Value = 0.111
$ws.Range("Y180").Value = -1
$ws.Range("Z180").Value = 0.875
$ws.Range("AA180").Value = -1
$ws.Range("AB180").Value = -0.5
$ws.Range("AC180").Value = 0.425

# Row 181
$ws.Range("B181").Value = 6937269
$ws.Range("F181").Value = "Aris Salonika"
$ws.Range("G181").Value = "AEK Athens"
$ws.Range("H181").Value = 3
$ws.Range("I181").Value = 3
$ws.Range("K181").Value = 4.75
$ws.Range("L181").Value = 3.75
$ws.Range("M181").Value = 1.75
$ws.Range("N181").Value = 6.5
$ws.Range("O181").Value = 4.2
$ws.Range("P181").Value = 1.5
$ws.Range("Q181").Value = 1
$ws.Range("R181").Value = 2.05
$ws.Range("S181").Value = 1.8
$ws.Range("T181").Value = 2.5
$ws.Range("U181").Value = 1.975
$ws.Range("V181").Value = 1.875
$ws.Range("X181").Value = 3.2
$ws.Range("Z181").Value = 1.05
$ws.Range("AA181").Value = -1
$ws.Range("AB181").Value = 0.9750000000000001
$ws.Range("AC181").Value = -1

# Row 189
$ws.Range("N189").Value = 2.2
$ws.Range("P189").Value = 3.6
$ws.Range("R189").Value = 1.9
$ws.Range("S189").Value = 1.95
$ws.Range("U189").Value = 2
$ws.Range("V189").Value = 1.85

# Row 190
$ws.Range("N190").Value = 2.25
$ws.Range("P190").Value = 3.3
$ws.Range("R190").Value = 1.95
$ws.Range("S190").Value = 1.9
$ws.Range("U190").Value = 1.9
$ws.Range("V190").Value = 1.95

# Row 191
$ws.Range("N191").Value = 2.05
$ws.Range("O191").Value = 3.25
$ws.Range("P191").Value = 4
$ws.Range("R191").Value = 2.05
$ws.Range("S191").Value = 1.8
$ws.Range("U191").Value = 1.975
$ws.Range("V191").Value = 1.875

# Row 192
$ws.Range("N192").Value = 2
$ws.Range("O192").Value = 3.5
$ws.Range("P192").Value = 3.75
$ws.Range("Q192").Value = -0.5
$ws.Range("R192").Value = 2.025
$ws.Range("S192").Value = 1.825
$ws.Range("U192").Value = 2
$ws.Range("V192").Value = 1.85

# Row 193
$ws.Range("P193").Value = 3.6
$ws.Range("R193").Value = 1.85
$ws.Range("S193").Value = 2

# Row 195
$ws.Range("N195").Value = 1.8
$ws.Range("O195").Value = 3.4
$ws.Range("P195").Value = 5
$ws.Range("R195").Value = 1.8
$ws.Range("S195").Value = 2.05
